# Added bulk email template
# Replace the author's personal placeholder addresses with generic
# template placeholders on both sheets, and nudge a couple of column
# widths / the remembered selection to match the refreshed template.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Sheet1
$ws2 = $wb.Worksheets.Item(2)   # Sheet2

# --- Sheet1: generic "reply to" email address used for the initial send ---
$ws1.Range("B2").Value = "EmailAddress@email.com"

# --- Sheet2: generic sender address for both the initial and reminder rows ---
$ws2.Range("C2").Value = "Sender@senderemail.com"
$ws2.Range("C3").Value = "Sender@senderemail.com"

# --- Column width tweaks on Sheet1 (B and D got a little wider) ---
$ws1.Columns.Item(2).ColumnWidth = 23.6667
$ws1.Columns.Item(4).ColumnWidth = 16.6667

# --- Restore/update the saved selections on both sheets ---
# (select Sheet2's range first so Sheet1 ends up as the active tab again)
$ws2.Range("C2:C3").Select()
$ws1.Range("B3").Select()

Write-Output "done"
